$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as text in the original workbook
# (e.g. "231.78", "37.158.75"). Excel would otherwise auto-convert plain
# decimal-looking strings into numbers, so force a Text number format on
# each Price cell we are about to update before assigning the new value.
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D9', 'D10', 'D11', 'D13', 'D14', 'D16', 'D17', 'D18', 'D19', 'D20', 'D22', 'D23', 'D25', 'D27', 'D28', 'D29', 'D30', 'D32', 'D33', 'D34', 'D35', 'D36', 'D42', 'D43', 'D44', 'D46', 'D47', 'D48', 'D50', 'D51')
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.129.74'
$ws.Range('E2').Value = '  +1.77%  '
$ws.Range('D3').Value = '2.049.47'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '231.78'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').Value = '0.618'
$ws.Range('E6').Value = '  +3.15%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +3.68%  '
$ws.Range('D9').Value = '0.381'
$ws.Range('E9').Value = '  +3.40%  '
$ws.Range('D10').Value = '57.43'
$ws.Range('D11').Value = '0.0755'
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('D13').Value = '2.355.03'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = '14.23'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('E15').Value = '  +4.07%  '
$ws.Range('D16').Value = '0.771'
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('D17').Value = '5.13'
$ws.Range('E17').Value = '  +0.96%  '
$ws.Range('D18').Value = '2.052.93'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').Value = '37.054.30'
$ws.Range('E19').Value = '  +1.62%  '
$ws.Range('D20').Value = '6.33'
$ws.Range('E20').Value = '  +8.38%  '
$ws.Range('E21').Value = '  +2.65%  '
$ws.Range('D22').Value = '0.0₃0806'
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('D23').Value = '224.90'
$ws.Range('E23').Value = '  +2.26%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').Value = '2.39'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').Value = '165.87'
$ws.Range('E27').Value = '  +1.99%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '8.74'
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').Value = '1.44'
$ws.Range('E29').Value = '  +6.54%  '
$ws.Range('D30').Value = '18.93'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('D32').Value = '0.117'
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').Value = '4.43'
$ws.Range('E33').Value = '  +2.07%  '
$ws.Range('D34').Value = '0.0614'
$ws.Range('E34').Value = '  +2.35%  '
$ws.Range('D35').Value = '4.56'
$ws.Range('E35').Value = '  +7.54%  '
$ws.Range('D36').Value = '2.50'
$ws.Range('E36').Value = '  +0.87%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').Value = '  -0.85%  '
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('D42').Value = '1.478.36'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '4.38'
$ws.Range('E43').Value = '  -3.63%  '
$ws.Range('D44').Value = '96.36'
$ws.Range('E44').Value = '  +3.48%  '
$ws.Range('E45').Value = '  +5.74%  '
$ws.Range('D46').Value = '0.0928'
$ws.Range('E46').Value = '  -1.04%  '
$ws.Range('D47').Value = '0.0209'
$ws.Range('E47').Value = '  +3.02%  '
$ws.Range('D48').Value = '1.01'
$ws.Range('E48').Value = '  +1.39%  '
$ws.Range('E49').Value = '  +3.57%  '
$ws.Range('D50').Value = '15.02'
$ws.Range('E50').Value = '  -3.01%  '
$ws.Range('D51').Value = '2.93'
$ws.Range('E51').Value = '  +1.35%  '

# Restore the default (unformatted) style on those cells now that the text
# value has been stored, so no visible/number-format change persists.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
